$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume data (Price = column D, Volume(1h) = column E)
$updates = @(
    @{ Row = 2; D = "29.223.21"; E = "  -0.52%  " }
    @{ Row = 3; D = "1.827.67"; E = "  -0.83%  " }
    @{ Row = 4; D = "1.005"; E = "  +0.58%  " }
    @{ Row = 5; D = "234.55"; E = "  -1.99%  " }
    @{ Row = 6; D = "0.5967"; E = "  -4.86%  " }
    @{ Row = 7; D = $null; E = "  +0.42%  " }
    @{ Row = 8; D = "0.06963"; E = "  -5.89%  " }
    @{ Row = 9; D = "0.2746"; E = "  -5.09%  " }
    @{ Row = 10; D = "23.23"; E = "  -6.52%  " }
    @{ Row = 11; D = "0.07620"; E = "  -1.20%  " }
    @{ Row = 12; D = "1.836.37"; E = "  -0.39%  " }
    @{ Row = 13; D = "4.756"; E = "  -4.28%  " }
    @{ Row = 14; D = "0.6241"; E = "  -7.48%  " }
    @{ Row = 15; D = "0.000009713"; E = "  -5.14%  " }
    @{ Row = 16; D = "78.39"; E = "  -4.30%  " }
    @{ Row = 17; D = "28.940.24"; E = "  -1.48%  " }
    @{ Row = 18; D = $null; E = "  -8.50%  " }
    @{ Row = 19; D = "221.80"; E = "  -5.37%  " }
    @{ Row = 20; D = $null; E = "  +0.40%  " }
    @{ Row = 21; D = "11.54"; E = "  -6.30%  " }
    @{ Row = 22; D = "6.878"; E = "  -5.90%  " }
    @{ Row = 23; D = "1.006"; E = "  +0.49%  " }
    @{ Row = 24; D = "155.94"; E = "  -1.12%  " }
    @{ Row = 25; D = "7.938"; E = "  -6.46%  " }
    @{ Row = 26; D = "0.1289"; E = "  -4.30%  " }
    @{ Row = 27; D = "16.48"; E = "  -4.85%  " }
    @{ Row = 28; D = "0.06663"; E = "  -8.58%  " }
    @{ Row = 29; D = "1.448"; E = "  -2.31%  " }
    @{ Row = 30; D = "1.439"; E = "  -2.60%  " }
    @{ Row = 31; D = "3.827"; E = "  -5.04%  " }
    @{ Row = 32; D = "3.751"; E = "  -6.97%  " }
    @{ Row = 33; D = "1.088"; E = "  -5.07%  " }
    @{ Row = 34; D = "1.712"; E = "  -5.82%  " }
    @{ Row = 35; D = "0.6410"; E = "  -8.40%  " }
    @{ Row = 36; D = "2.549"; E = "  -0.83%  " }
    @{ Row = 37; D = "2.731"; E = "  -2.33%  " }
    @{ Row = 38; D = "1.186.58"; E = "  -3.83%  " }
    @{ Row = 39; D = "0.01732"; E = "  -5.61%  " }
    @{ Row = 40; D = "6.497"; E = "  -5.58%  " }
    @{ Row = 41; D = "0.9015"; E = "  -4.91%  " }
    @{ Row = 42; D = $null; E = "  +0.49%  " }
    @{ Row = 43; D = "1.980.32"; E = "  -0.53%  " }
    @{ Row = 44; D = "100.29"; E = "  -0.68%  " }
    @{ Row = 45; D = "61.93"; E = "  -5.09%  " }
    @{ Row = 46; D = "0.00000000114"; E = "  -4.79%  " }
    @{ Row = 47; D = "8.436"; E = "  -4.86%  " }
    @{ Row = 48; D = "0.4559"; E = "  -0.26%  " }
    @{ Row = 49; D = "0.05508"; E = "  -2.68%  " }
    @{ Row = 50; D = "1.568"; E = "  -8.12%  " }
    @{ Row = 51; D = "6.341"; E = "  -8.96%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }

    $cellE = $ws.Cells.Item($u.Row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
    $cellE.Style = "Normal"
}
